$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 239, shifting existing rows 239-258 down to 241-260.
$ws.Rows.Item(239).Insert()
$ws.Rows.Item(239).Insert()

# Fill in new row 239 with its data.
$ws.Cells.Item(239, 1).Value = 7
$ws.Cells.Item(239, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(239, 3).Value = "Ñuble"
$ws.Cells.Item(239, 4).Value = 44783
$ws.Cells.Item(239, 5).Value = 16
$ws.Cells.Item(239, 6).Value = 100112009
$ws.Cells.Item(239, 7).Value = "Acelga"
$ws.Cells.Item(239, 8).Value = "Sin especificar"
$ws.Cells.Item(239, 9).Value = "Primera"
$ws.Cells.Item(239, 10).Value = 200
$ws.Cells.Item(239, 11).Value = 700
$ws.Cells.Item(239, 12).Value = 800
$ws.Cells.Item(239, 13).Value = 750
$ws.Cells.Item(239, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(239, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(239, 16).Value = 750
$ws.Cells.Item(239, 17).Value = 1
$ws.Cells.Item(239, 18).Value = "Hortaliza"

# Fill in new row 240 with its data.
$ws.Cells.Item(240, 1).Value = 7
$ws.Cells.Item(240, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(240, 3).Value = "Ñuble"
$ws.Cells.Item(240, 4).Value = 44783
$ws.Cells.Item(240, 5).Value = 16
$ws.Cells.Item(240, 6).Value = 100112009
$ws.Cells.Item(240, 7).Value = "Acelga"
$ws.Cells.Item(240, 8).Value = "Sin especificar"
$ws.Cells.Item(240, 9).Value = "Segunda"
$ws.Cells.Item(240, 10).Value = 150
$ws.Cells.Item(240, 11).Value = 600
$ws.Cells.Item(240, 12).Value = 600
$ws.Cells.Item(240, 13).Value = 600
$ws.Cells.Item(240, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(240, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(240, 16).Value = 600
$ws.Cells.Item(240, 17).Value = 1
$ws.Cells.Item(240, 18).Value = "Hortaliza"
